# Refresh the "cryptos" price/volume table (Price = col D, Volume(1h) = col E)
# for rows 2-51 with the latest scraped figures.
#
# Note: a handful of "Price" values (e.g. "231.54", "9.70") look like plain
# decimals to Excel's auto-detection, so a bare .Value assignment would turn
# them into numbers (losing significant trailing zeros, e.g. "9.70"->"9.7").
# Prefixing with a leading apostrophe forces them to stay text (matching the
# source data, which stores every Price/Volume cell as text); ClearFormats()
# afterwards drops the "Text" number-format Excel applies along with that
# apostrophe, so the cell keeps its original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.079.31"
$ws.Range("E2").Value = "  +5.80%  "
$ws.Range("D3").Value = "2.234.57"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'231.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D7").Value = "'61.81"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.23%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.84%  "
$ws.Range("D10").Value = "'58.58"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").Value = "2.568.23"
$ws.Range("E13").Value = "  +2.49%  "
$ws.Range("E14").Value = "  -2.01%  "
$ws.Range("D15").Value = "'22.01"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'0.804"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "2.246.57"
$ws.Range("E18").Value = "  +3.61%  "
$ws.Range("D19").Value = "41.929.84"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").Value = "'72.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "0.0₃0898"
$ws.Range("E21").Value = "  -2.84%  "
$ws.Range("D22").Value = "'6.03"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "'250.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +8.65%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.39"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("E26").Value = "  -3.42%  "
$ws.Range("D27").Value = "'9.70"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("D29").Value = "'167.59"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.92%  "
$ws.Range("D30").Value = "'20.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E32").Value = "  +1.28%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'5.04"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.01%  "
$ws.Range("D35").Value = "'4.66"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'6.63"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.26%  "
$ws.Range("D38").Value = "'3.70"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.30%  "
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("D40").Value = "'0.000269"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +39.27%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "'0.0240"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.25%  "
$ws.Range("D43").Value = "'4.83"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.35%  "
$ws.Range("D44").Value = "'8.53"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.69%  "
$ws.Range("E45").Value = "  +5.68%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'99.01"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.79%  "
$ws.Range("D48").Value = "1.479.84"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("E49").Value = "  -7.10%  "
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "'52.70"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.32%  "
